# Apply odds updates per commit "Atualizando o arquivo XLSX"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("J3").Value = 2.87
$ws.Range("K3").Value = 1.92
$ws.Range("O3").Value = 1.5
$ws.Range("P3").Value = 2.5
$ws.Range("Q3").Value = 1.9
$ws.Range("R3").Value = 1.95
$ws.Range("S3").Value = 2.5
$ws.Range("T3").Value = 1.5
$ws.Range("U3").Value = 3.95
# Row 4
$ws.Range("K4").Value = 1.87
# Row 5
$ws.Range("G5").Value = 3.15
$ws.Range("H5").Value = 2.45
$ws.Range("I5").Value = 2.85
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 1.7
$ws.Range("M5").Value = 1.19
$ws.Range("N5").Value = 4.15
$ws.Range("O5").Value = 1.75
$ws.Range("P5").Value = 1.98
$ws.Range("S5").Value = 3.15
$ws.Range("T5").Value = 1.31
$ws.Range("W5").Value = 5.9
$ws.Range("X5").Value = 1.1
$ws.Range("Z5").Value = 2
$ws.Range("AA5").Value = 2.32
$ws.Range("AB5").Value = 1.53
$ws.Range("AC5").Value = 6
$ws.Range("AE5").Value = 12.5
$ws.Range("AH5").Value = 70
$ws.Range("AI5").Value = 4.15
$ws.Range("AK5").Value = 20
$ws.Range("AL5").Value = 175
$ws.Range("AN5").Value = 5.7
$ws.Range("AP5").Value = 11.5
$ws.Range("AS5").Value = 65
# Row 6
$ws.Range("G6").Value = 2.85
$ws.Range("H6").Value = 2.65
$ws.Range("I6").Value = 2.82
$ws.Range("K6").Value = 1.83
$ws.Range("L6").Value = 3.55
$ws.Range("M6").Value = 1.16
$ws.Range("N6").Value = 4.65
$ws.Range("O6").Value = 1.65
$ws.Range("P6").Value = 2.12
$ws.Range("S6").Value = 2.87
$ws.Range("T6").Value = 1.36
$ws.Range("W6").Value = 5.3
$ws.Range("X6").Value = 1.12
$ws.Range("Y6").Value = 1.6
$ws.Range("Z6").Value = 2.2
$ws.Range("AA6").Value = 2.25
$ws.Range("AB6").Value = 1.57
$ws.Range("AC6").Value = 6.1
$ws.Range("AD6").Value = 12.5
$ws.Range("AE6").Value = 11.25
$ws.Range("AH6").Value = 55
$ws.Range("AI6").Value = 4.65
$ws.Range("AK6").Value = 19.5
$ws.Range("AL6").Value = 150
$ws.Range("AN6").Value = 5.9
$ws.Range("AO6").Value = 12
$ws.Range("AP6").Value = 11.5
$ws.Range("AQ6").Value = 37
$ws.Range("AR6").Value = 35
$ws.Range("AS6").Value = 60
# Row 7
$ws.Range("H7").Value = 2.9
$ws.Range("I7").Value = 3.5
$ws.Range("J7").Value = 3.25
$ws.Range("K7").Value = 1.8
$ws.Range("L7").Value = 4.5
$ws.Range("M7").Value = 1.14
$ws.Range("N7").Value = 5.5
$ws.Range("O7").Value = 1.67
$ws.Range("P7").Value = 2.2
$ws.Range("S7").Value = 3.2
$ws.Range("T7").Value = 1.36
$ws.Range("Y7").Value = 1.67
$ws.Range("Z7").Value = 2.1
$ws.Range("AA7").Value = 2.38
$ws.Range("AB7").Value = 1.53
$ws.Range("AD7").Value = 9
$ws.Range("AE7").Value = 11
$ws.Range("AG7").Value = 26
$ws.Range("AI7").Value = 5.5
$ws.Range("AL7").Value = 101
$ws.Range("AP7").Value = 15
$ws.Range("AR7").Value = 41
# Row 8
$ws.Range("G8").Value = 2.6
$ws.Range("I8").Value = 2.6
$ws.Range("J8").Value = 3.4
$ws.Range("L8").Value = 3.4
$ws.Range("M8").Value = 1.07
$ws.Range("N8").Value = 8.5
$ws.Range("O8").Value = 1.36
$ws.Range("P8").Value = 3
$ws.Range("S8").Value = 2.15
$ws.Range("T8").Value = 1.67
$ws.Range("W8").Value = 4
$ws.Range("X8").Value = 1.22
$ws.Range("AA8").Value = 1.91
$ws.Range("AB8").Value = 1.91
$ws.Range("AF8").Value = 26
$ws.Range("AG8").Value = 23
$ws.Range("AI8").Value = 8.5
$ws.Range("AN8").Value = 8
$ws.Range("AO8").Value = 12
$ws.Range("AP8").Value = 10
$ws.Range("AQ8").Value = 26
# Row 9
$ws.Range("G9").Value = 2.38
$ws.Range("I9").Value = 2.8
$ws.Range("J9").Value = 3.2
$ws.Range("L9").Value = 3.6
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 9
$ws.Range("AC9").Value = 7.5
$ws.Range("AD9").Value = 11
$ws.Range("AE9").Value = 9.5
$ws.Range("AF9").Value = 23
$ws.Range("AG9").Value = 21
$ws.Range("AK9").Value = 17
$ws.Range("AO9").Value = 13
$ws.Range("AP9").Value = 11
$ws.Range("AQ9").Value = 29
$ws.Range("AR9").Value = 23
# Row 10
$ws.Range("H10").Value = 3.3
$ws.Range("I10").Value = 2
$ws.Range("K10").Value = 2
$ws.Range("M10").Value = 1.08
$ws.Range("N10").Value = 7.5
$ws.Range("O10").Value = 1.44
$ws.Range("P10").Value = 2.63
$ws.Range("Q10").Value = 1.78
$ws.Range("R10").Value = 2.1
$ws.Range("S10").Value = 2.35
$ws.Range("T10").Value = 1.57
$ws.Range("W10").Value = 4.5
$ws.Range("X10").Value = 1.18
$ws.Range("Y10").Value = 1.53
$ws.Range("Z10").Value = 2.38
$ws.Range("AA10").Value = 2.05
$ws.Range("AB10").Value = 1.7
$ws.Range("AC10").Value = 8.5
$ws.Range("AD10").Value = 17
$ws.Range("AI10").Value = 7.5
$ws.Range("AK10").Value = 19
$ws.Range("AL10").Value = 67
$ws.Range("AM10").Value = 501
$ws.Range("AN10").Value = 6
$ws.Range("AP10").Value = 9.5
$ws.Range("AR10").Value = 19
$ws.Range("AS10").Value = 34
# Row 11
$ws.Range("H11").Value = 3.3
$ws.Range("K11").Value = 2.05
$ws.Range("Y11").Value = 1.5
$ws.Range("Z11").Value = 2.5
$ws.Range("AA11").Value = 2.05
$ws.Range("AB11").Value = 1.7
$ws.Range("AC11").Value = 6
$ws.Range("AD11").Value = 7.5
$ws.Range("AE11").Value = 9
$ws.Range("AG11").Value = 17
$ws.Range("AH11").Value = 34
$ws.Range("AI11").Value = 8
$ws.Range("AK11").Value = 19
$ws.Range("AL11").Value = 67
$ws.Range("AM11").Value = 451
$ws.Range("AO11").Value = 21
$ws.Range("AS11").Value = 51
# Row 14
$ws.Range("G14").Value = 2.05
$ws.Range("I14").Value = 3.4
$ws.Range("J14").Value = 2.63
$ws.Range("M14").Value = 1.03
$ws.Range("N14").Value = 9.5
$ws.Range("AB14").Value = 1.91
$ws.Range("AD14").Value = 10
$ws.Range("AF14").Value = 19
$ws.Range("AN14").Value = 10
$ws.Range("AO14").Value = 17
